$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: populate the soon-to-be G/H/I columns while the sheet still
# only has columns A-E, writing to F/G/H so that after the later
# column-insert at C they land in G/H/I (matches shared-string append
# order captured in the target sharedStrings.xml).
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "slrtype"
$ws.Range("F1").Style = "Normal"
$ws.Range("G1").Value = "slrtype_Radio_button"
$ws.Range("G1").Style = "Normal"
$ws.Range("H1").Value = "data"

$ws.Range("F2").Value = "Clinical"
$ws.Range("G2").Value = "Clinical_radio_button"
$ws.Range("H2").Value = "Adult patients (18 years or older) with maintenance therapy for NDMM after induction therapy, post-SCT*"

$ws.Range("F3").Value = "Economic"
$ws.Range("G3").Value = "Economic_radio_button"
$ws.Range("H3").Value = "Adult patients (18 years or older) with maintenance therapy for NDMM after induction therapy, post-SCT*"

$ws.Range("F4").Value = "Quality of Life"
$ws.Range("G4").Value = "Quality of Life_radio_button"
$ws.Range("H4").Value = "Patients who do not undergo maintenance therapy`nPatients in later lines of therapy (not first line)`nPediatric population`nOther types of cancers"
$ws.Range("H4").WrapText = $true
$ws.Rows("4").RowHeight = 100.8

$ws.Range("F5").Value = "Real-world Evidence"
$ws.Range("G5").Value = "Real-world Evidence_radio_button"

# ---------------------------------------------------------------------
# Step 2: insert the new column C (pushes old C/D/E -> D/E/F, and the
# F/G/H we just filled -> G/H/I).
# ---------------------------------------------------------------------
$ws.Columns("C").Insert(-4161)
$ws.Range("C1").Value = "Population_Radio_button"
$ws.Range("C1").Style = "Normal"
$ws.Range("C2").Value = "Test - Test_radio_button"

# New column C should share column B's width and have no bestFit.
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# ---------------------------------------------------------------------
# Match the author's final selection (cell C4).
# ---------------------------------------------------------------------
$ws.Range("C4").Select()
